# Bugfixes commit: adds a new "Form" worksheet (translation table for
# frm_Settings / frm_editFiledata) between "FileBrowser" and "GroupBox",
# makes it the active sheet, and tweaks a couple of neighbouring sheet
# selections (FileBrowser's selection becomes a full-range select, and the
# tabSelected flag moves off "ColumnHeader" onto the new "Form" sheet
# automatically because it becomes the active one).

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "Form" sheet right after "FileBrowser" ----------
$fileBrowserSheet = $wb.Worksheets.Item("FileBrowser")
$groupBoxSheet    = $wb.Worksheets.Item("GroupBox")

$formSheet = $wb.Worksheets.Add($null, $fileBrowserSheet)
$formSheet.Name = "Form"

# --- 2. Populate header row + data, matching the other lookup sheets ---
$formSheet.Range("A1").Value = "objectName"
$formSheet.Range("B1").Value = "actionType"
$formSheet.Range("C1").Value = "objectText"

$formSheet.Range("A2").Value = "frm_Settings"
$formSheet.Range("C2").Value = "Settings"

$formSheet.Range("A3").Value = "frm_editFiledata"
$formSheet.Range("C3").Value = "Edit Data"

# Reuse the exact header style (bold/indexed-black font, xf index 4) that
# every other lookup sheet already uses, by copying formats from GroupBox
# instead of re-deriving a (slightly different) style via Font.Bold.
$groupBoxSheet.Range("A1:C1").Copy() | Out-Null
$formSheet.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- 3. Column widths (best effort match of the authored widths) -------
$formSheet.Columns.Item(1).ColumnWidth = 11.3046875
$formSheet.Columns.Item(2).ColumnWidth = 10
$formSheet.Columns.Item(3).ColumnWidth = 9.69140625

# --- 4. Freeze the header row and select C4, make Form the active tab --
$formSheet.Activate()
$formSheet.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$formSheet.Range("C4").Select()

# --- 5. FileBrowser: selection becomes a full A1:C2 range select -------
$fileBrowserSheet.Activate()
$fileBrowserSheet.Range("A1:C2").Select()

# Leave "Form" as the active sheet/tab, matching the saved workbook state.
$formSheet.Activate()
